$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 3959
$ws.Range("J17").Value = 3959
$ws.Range("L17").Value = 11877
$ws.Range("N17").Value = -12213

$ws.Range("H98").Value = 1190
$ws.Range("I98").Value = 650
$ws.Range("J98").Value = 2000
$ws.Range("K98").Value = 650
$ws.Range("L98").Value = 2000
$ws.Range("M98").Value = 848
$ws.Range("N98").Value = -4996

$ws.Range("H100").Value = 835102.7
$ws.Range("I100").Value = 2002194.4
$ws.Range("J100").Value = 1465.7142
$ws.Range("K100").Value = 2002194.4
$ws.Range("L100").Value = 1465.7142
$ws.Range("M100").Value = -2001653.4
$ws.Range("N100").Value = -2547.7142

$ws.Range("H112").Value = 1233.2963
$ws.Range("J112").Value = 1892.7858
$ws.Range("L112").Value = 5678.357400000001
$ws.Range("N112").Value = -7894.357400000001

$ws.Range("H113").Value = 3335333
$ws.Range("I113").Value = 5001499.5
$ws.Range("J113").Value = 3000
$ws.Range("K113").Value = 5001499.5
$ws.Range("L113").Value = 3000
$ws.Range("M113").Value = -4998245.5
$ws.Range("N113").Value = -9508

$ws.Range("H122").Value = 1190
$ws.Range("I122").Value = 650
$ws.Range("J122").Value = 2000
$ws.Range("K122").Value = 1950
$ws.Range("L122").Value = 6000
$ws.Range("M122").Value = 500
$ws.Range("N122").Value = -10900

$ws.Range("H138").Value = 2739.1091
$ws.Range("I138").Value = 1818.8462
$ws.Range("J138").Value = 3023.9524
$ws.Range("K138").Value = 5456.5386
$ws.Range("L138").Value = 9071.8572
$ws.Range("M138").Value = -316.5385999999999
$ws.Range("N138").Value = -19351.8572

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 973.4
$ws.Range("I2").Value = 1108.75
$ws.Range("J2").Value = 432
$ws.Range("K2").Value = 1108.75
$ws.Range("L2").Value = 432
$ws.Range("M2").Value = -995.75
$ws.Range("N2").Value = -658

$ws.Range("H61").Value = 2463.72
$ws.Range("I61").Value = 1892.9231
$ws.Range("K61").Value = 1892.9231
$ws.Range("M61").Value = -1680.9231

$ws.Range("H74").Value = 39986524
$ws.Range("J74").Value = 2997.5
$ws.Range("L74").Value = 2997.5
$ws.Range("N74").Value = -4745.5

$ws.Range("H77").Value = 39986524
$ws.Range("J77").Value = 2997.5
$ws.Range("L77").Value = 14987.5
$ws.Range("N77").Value = -23723.5

$ws.Range("H98").Value = 70000
$ws.Range("J98").Value = 70000
$ws.Range("L98").Value = 70000
$ws.Range("N98").Value = -75990

$ws.Range("H116").Value = 973.4
$ws.Range("I116").Value = 1108.75
$ws.Range("J116").Value = 432
$ws.Range("K116").Value = 1108.75
$ws.Range("L116").Value = 432
$ws.Range("M116").Value = 1185.25
$ws.Range("N116").Value = -5020

$ws.Range("H136").Value = 2463.72
$ws.Range("I136").Value = 1892.9231
$ws.Range("K136").Value = 5678.7693
$ws.Range("M136").Value = -3128.7693

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 973.4
$ws.Range("I3").Value = 1108.75
$ws.Range("J3").Value = 432
$ws.Range("K3").Value = 1108.75
$ws.Range("L3").Value = 432
$ws.Range("M3").Value = -994.75
$ws.Range("N3").Value = -660

$ws.Range("H26").Value = 50000
$ws.Range("I26").Value = 45000
$ws.Range("K26").Value = 45000
$ws.Range("M26").Value = -44708

$ws.Range("H28").Value = 99500
$ws.Range("J28").Value = 99500
$ws.Range("L28").Value = 99500
$ws.Range("N28").Value = -100088

$ws.Range("H87").Value = 89999
$ws.Range("J87").Value = 89999
$ws.Range("L87").Value = 89999
$ws.Range("N87").Value = -92495

$ws.Range("H90").Value = 89999
$ws.Range("J90").Value = 89999
$ws.Range("L90").Value = 269997
$ws.Range("N90").Value = -282477

$ws.Range("H94").Value = 502.7143
$ws.Range("I94").Value = 502.7143
$ws.Range("K94").Value = 502.7143
$ws.Range("M94").Value = -51.71429999999998

$ws.Range("H105").Value = 4578.8335
$ws.Range("I105").Value = 4618.5
$ws.Range("K105").Value = 4618.5
$ws.Range("M105").Value = -2871.5

$ws.Range("H134").Value = 2178.6667
$ws.Range("I134").Value = 1991.2941
$ws.Range("K134").Value = 5973.8823
$ws.Range("M134").Value = -3438.8823

$ws.Range("H140").Value = 74999
$ws.Range("I140").Value = 74999
$ws.Range("J140").Value = 0
$ws.Range("K140").Value = 74999
$ws.Range("L140").Value = 0
$ws.Range("M140").Value = -69819
$ws.Range("N140").Value = ""

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 377.7143
$ws.Range("I7").Value = 258.8
$ws.Range("K7").Value = 258.8
$ws.Range("M7").Value = -145.8

$ws.Range("H16").Value = 1283.125
$ws.Range("I16").Value = 1378.8334
$ws.Range("K16").Value = 1378.8334
$ws.Range("M16").Value = -1091.8334

$ws.Range("H31").Value = 4133
$ws.Range("I31").Value = 4133
$ws.Range("K31").Value = 4133
$ws.Range("M31").Value = -3838

$ws.Range("H34").Value = 4133
$ws.Range("I34").Value = 4133
$ws.Range("K34").Value = 4133
$ws.Range("M34").Value = -3931

$ws.Range("H94").Value = 674.75
$ws.Range("J94").Value = 233
$ws.Range("L94").Value = 233
$ws.Range("N94").Value = -1135

$ws.Range("H113").Value = 1283.125
$ws.Range("I113").Value = 1378.8334
$ws.Range("K113").Value = 1378.8334
$ws.Range("M113").Value = 791.1666

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 1078.5714
$ws.Range("I131").Value = 710.6
$ws.Range("J131").Value = 1998.5
$ws.Range("K131").Value = 2131.8
$ws.Range("L131").Value = 5995.5
$ws.Range("M131").Value = 2908.2
$ws.Range("N131").Value = -16075.5

$ws.Range("H140").Value = 1254.5
$ws.Range("I140").Value = 1254.5
$ws.Range("K140").Value = 3763.5
$ws.Range("M140").Value = 1416.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 917.1429000000001
$ws.Range("I97").Value = 736.6667
$ws.Range("J97").Value = 2000
$ws.Range("K97").Value = 736.6667
$ws.Range("L97").Value = 2000
$ws.Range("M97").Value = -240.6667
$ws.Range("N97").Value = -2992

$ws.Range("H122").Value = 1998.5
$ws.Range("I122").Value = 1999
$ws.Range("J122").Value = 1998
$ws.Range("K122").Value = 5997
$ws.Range("L122").Value = 5994
$ws.Range("M122").Value = -3547
$ws.Range("N122").Value = -10894

$ws.Range("H132").Value = 2783.353
$ws.Range("I132").Value = 1702.5555
$ws.Range("J132").Value = 3999.25
$ws.Range("K132").Value = 5107.666499999999
$ws.Range("L132").Value = 11997.75
$ws.Range("M132").Value = -2577.666499999999
$ws.Range("N132").Value = -17057.75

$ws.Range("H134").Value = 110163
$ws.Range("J134").Value = 110163
$ws.Range("L134").Value = 330489
$ws.Range("N134").Value = -335559

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 4983.1665
$ws.Range("I40").Value = 4983.1665
$ws.Range("J40").Value = 0
$ws.Range("K40").Value = 4983.1665
$ws.Range("L40").Value = 0
$ws.Range("M40").Value = -4847.1665
$ws.Range("N40").Value = ""

$ws.Range("H82").Value = 1866
$ws.Range("I82").Value = 1798
$ws.Range("J82").Value = 1900
$ws.Range("K82").Value = 1798
$ws.Range("L82").Value = 1900
$ws.Range("M82").Value = -1437
$ws.Range("N82").Value = -2622

$ws.Range("H85").Value = 1866
$ws.Range("I85").Value = 1798
$ws.Range("J85").Value = 1900
$ws.Range("K85").Value = 1798
$ws.Range("L85").Value = 1900
$ws.Range("M85").Value = -550
$ws.Range("N85").Value = -4396

$ws.Range("H93").Value = 3517
$ws.Range("I93").Value = 4020.4
$ws.Range("J93").Value = 1000
$ws.Range("K93").Value = 4020.4
$ws.Range("L93").Value = 1000
$ws.Range("M93").Value = -2772.4
$ws.Range("N93").Value = -3496

$ws.Range("H136").Value = 13334816
$ws.Range("I136").Value = 13334816
$ws.Range("K136").Value = 40004448
$ws.Range("M136").Value = -40001898

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H51").Value = 0
$ws.Range("J51").Value = 0
$ws.Range("L51").Value = 0
$ws.Range("N51").Value = ""

$ws.Range("H70").Value = 52323.332
$ws.Range("I70").Value = 45000
$ws.Range("J70").Value = 53788
$ws.Range("K70").Value = 45000
$ws.Range("L70").Value = 53788
$ws.Range("M70").Value = -44685
$ws.Range("N70").Value = -54418

$ws.Range("H73").Value = 52323.332
$ws.Range("I73").Value = 45000
$ws.Range("J73").Value = 53788
$ws.Range("K73").Value = 45000
$ws.Range("L73").Value = 53788
$ws.Range("M73").Value = -43908
$ws.Range("N73").Value = -55972

$ws.Range("H81").Value = 12496.6
$ws.Range("I81").Value = 8166.6665
$ws.Range("J81").Value = 18991.5
$ws.Range("K81").Value = 16333.333
$ws.Range("L81").Value = 37983
$ws.Range("M81").Value = -15272.333
$ws.Range("N81").Value = -40105

$ws.Range("H84").Value = 12496.6
$ws.Range("I84").Value = 8166.6665
$ws.Range("J84").Value = 18991.5
$ws.Range("K84").Value = 81666.66500000001
$ws.Range("L84").Value = 189915
$ws.Range("M84").Value = -76362.66500000001
$ws.Range("N84").Value = -200523

$ws.Range("H100").Value = 1385
$ws.Range("I100").Value = 1139.5
$ws.Range("J100").Value = 1876
$ws.Range("K100").Value = 2279
$ws.Range("L100").Value = 3752
$ws.Range("M100").Value = -1738
$ws.Range("N100").Value = -4834
